$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "vhbh jfv j"
$ws.Range("B1").Value = "jkbkvkh"
$ws.Range("C1").Value = "v"
$ws.Range("D1").Value = "hkvk"
$ws.Range("E1").Value = "hh"

$ws.Range("E1").Select() | Out-Null
